{"js": "// Stabilize app without exposing API keys:\n// Append two new paragraphs at the end of the document body \u2014 a\n// \"newapi\" label paragraph followed by a paragraph holding the new\n// (replacement) API key \u2014 mirroring the existing \"api=\" / key pair.\nconst body = context.document.body;\n\nbody.insertParagraph(\"newapi\", Word.InsertLocation.end);\nbody.insertParagraph(\n  \"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\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Stabilize app without exposing API keys:\n# Append two new paragraphs at the end of the document - a \"newapi\"\n# label paragraph followed by a paragraph holding the new (replacement)\n# API key - mirroring the existing \"api=\" / key pair.\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Last\n$p1.Range.Text = \"newapi\"\n\n$r2 = $d.Content\n$r2.Collapse(0)\n$r2.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = \"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\"\n"}
